$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "culture_collection" field/column (AH) is being removed from the MIGS
# template again (re-deletion per INSDC2017 review). Deleting the whole
# column shifts every later column (data + header shared-strings) one slot
# to the left, which matches the diff for sheet1.xml and sharedStrings.xml.
$ws.Columns("AH").Delete()

# Cell comments are NOT automatically moved by the column delete on this
# runtime, so re-home each remaining header comment (row 15) by hand: every
# comment from AH15..BH15 must show the text that used to belong to the
# next column over (AI15..BI15), and the now-nonexistent BI15 comment must
# be removed entirely.
$cols = @("AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY","AZ","BA","BB","BC","BD","BE","BF","BG","BH")

$newTexts = @(
    "temperature to which a given parcel of humid air must be cooled, at constant barometric pressure, for water vapor to condense into water.",
    "Traits like antibiotic resistance/xenobiotic degration phenotypes/converting phage genes",
    "Estimated size of genome",
    "Plasmids that have significance phenotypic consequence",
    "Health or disease status of sample at time of collection",
    "The natural (as opposed to laboratory) host to the organism from which the sample was obtained. Use the full taxonomic name, eg, ""Homo sapiens"".",
    "NCBI taxonomy ID of the host, e.g. 9606",
    "type of indoor surface",
    "Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.",
    "A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html",
    "To what is the entity pathogenic",
    "Aerobic or anaerobic",
    "Method or device employed for collecting sample",
    "Processing applied to the sample during or after isolation",
    "Amount or size of sample (volume, mass or area) that was collected",
    "method by which samples are sorted",
    "volume (mL) or weight (g) of sample processed for DNA extraction",
    "unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.",
    "Information about the genetic distinctness of the lineage (eg., biovar, serovar)",
    "substructure or under building is that largely hidden section of the building which is built off the foundations to the ground floor level",
    "contaminant identified on surface",
    "surfaces: water activity as a function of air and material moisture",
    "surface materials at the point of sampling",
    "water held on a surface",
    "pH measurement of surface",
    "temperature of the surface at the time of sampling",
    "Feeding position in food chain (eg., chemolithotroph)"
)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "15").Comment.Text($newTexts[$i]) | Out-Null
}

# BI no longer exists as a data column, so its trailing comment goes away.
$ws.Range("BI15").Comment.Delete() | Out-Null
